$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain number string (e.g. "1.002") would be
# auto-converted to a numeric type by Excel on assignment; temporarily mark
# them as Text so they are stored as strings, matching the source data, then
# clear the format again so no stray number-format is left behind.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.809.63'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '1.743.72'
$ws.Range('E3').Value = '  -2.12%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '333.36'
$ws.Range('E5').Value = '  -0.92%  '
$ws.Range('D6').Value = '0.9999'
$ws.Range('E6').Value = '  +0.50%  '
$ws.Range('D7').Value = '0.3883'
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('D8').Value = '0.3372'
$ws.Range('E8').Value = '  -1.62%  '
$ws.Range('D9').Value = '45.35'
$ws.Range('E9').Value = '  -4.21%  '
$ws.Range('D10').Value = '1.098'
$ws.Range('E10').Value = '  -4.89%  '
$ws.Range('D11').Value = '0.07146'
$ws.Range('E11').Value = '  -3.49%  '
$ws.Range('D12').Value = '0.9990'
$ws.Range('E12').Value = '  -0.19%  '
$ws.Range('D13').Value = '21.86'
$ws.Range('E13').Value = '  -5.51%  '
$ws.Range('D14').Value = '6.074'
$ws.Range('E14').Value = '  -5.38%  '
$ws.Range('D15').Value = '1.742.16'
$ws.Range('E15').Value = '  -2.59%  '
$ws.Range('D16').Value = '6.931'
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').Value = '0.00001047'
$ws.Range('E17').Value = '  -3.09%  '
$ws.Range('D18').Value = '0.06606'
$ws.Range('E18').Value = '  -0.72%  '
$ws.Range('D19').Value = '79.02'
$ws.Range('E19').Value = '  -4.47%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').Value = '16.73'
$ws.Range('E21').Value = '  -4.24%  '
$ws.Range('D22').Value = '6.157'
$ws.Range('E22').Value = '  -4.17%  '
$ws.Range('D23').Value = '27.778.27'
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').Value = '11.49'
$ws.Range('E24').Value = '  -4.96%  '
$ws.Range('D25').Value = '2.383'
$ws.Range('E25').Value = '  +0.45%  '
$ws.Range('D26').Value = '153.80'
$ws.Range('E26').Value = '  -0.43%  '
$ws.Range('D27').Value = '19.74'
$ws.Range('E27').Value = '  -5.39%  '
$ws.Range('D28').Value = '2.278'
$ws.Range('E28').Value = '  -5.35%  '
$ws.Range('D29').Value = '1.940.15'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('D30').Value = '1.271'
$ws.Range('E30').Value = '  -10.88%  '
$ws.Range('D31').Value = '127.48'
$ws.Range('E31').Value = '  -5.85%  '
$ws.Range('D32').Value = '4.049'
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').Value = '5.738'
$ws.Range('E33').Value = '  -6.74%  '
$ws.Range('D34').Value = '0.08686'
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('D35').Value = '11.94'
$ws.Range('E35').Value = '  -6.45%  '
$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').Value = '5.091'
$ws.Range('E36').Value = '  -4.52%  '
$ws.Range('B37').Value = 'WEMIXTOKEN'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').Value = '1.507'
$ws.Range('E37').Value = '  +0.10%  '
$ws.Range('D38').Value = '0.02254'
$ws.Range('E38').Value = '  -6.91%  '
$ws.Range('D39').Value = '0.06062'
$ws.Range('E39').Value = '  -4.60%  '
$ws.Range('D40').Value = '0.6405'
$ws.Range('E40').Value = '  -6.72%  '
$ws.Range('D41').Value = '0.2080'
$ws.Range('E41').Value = '  -4.45%  '
$ws.Range('D42').Value = '1.188'
$ws.Range('E42').Value = '  -4.17%  '
$ws.Range('D43').Value = '0.9994'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('D44').Value = '7.853'
$ws.Range('E44').Value = '  -5.49%  '
$ws.Range('D45').Value = '13.55'
$ws.Range('E45').Value = '  -5.78%  '
$ws.Range('D46').Value = '3.803'
$ws.Range('E46').Value = '  -1.28%  '
$ws.Range('D47').Value = '0.5914'
$ws.Range('E47').Value = '  -6.30%  '
$ws.Range('D48').Value = '125.73'
$ws.Range('E48').Value = '  -4.97%  '
$ws.Range('D49').Value = '1.966'
$ws.Range('E49').Value = '  -5.94%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').Value = '1.142'
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06913'
$ws.Range('E51').Value = '  -7.81%  '

$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()

